$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8147768276045895
$ws.Range("C2").Value = 0.2136217747381863
$ws.Range("D2").Value = 0.07946615628169695
$ws.Range("E2").Value = 0.1191513367428243
$ws.Range("G2").Value = 0.310941360334084
$ws.Range("H2").Value = 0.4658305808968066
$ws.Range("M2").Value = 0.3380024614860631
$ws.Range("O2").Value = 1.482401996685581
$ws.Range("B3").Value = 0.7131801136489457
$ws.Range("C3").Value = 0.1910689596984696
$ws.Range("D3").Value = 0.0719541466638276
$ws.Range("E3").Value = 0.1143489186497888
$ws.Range("G3").Value = 0.308828250439781
$ws.Range("H3").Value = 0.4696116846105838
$ws.Range("M3").Value = 0.2999044575238798
$ws.Range("O3").Value = 1.485516303262898
$ws.Range("B4").Value = 0.6506230520072052
$ws.Range("C4").Value = 0.1771560867453843
$ws.Range("D4").Value = 0.06737620766929808
$ws.Range("E4").Value = 0.1115336791581427
$ws.Range("G4").Value = 0.3079197556072728
$ws.Range("H4").Value = 0.4722743592603038
$ws.Range("M4").Value = 0.2765429980186482
$ws.Range("O4").Value = 1.488955868022387
$ws.Range("B5").Value = 0.6250875958371012
$ws.Range("C5").Value = 0.17147041271636
$ws.Range("D5").Value = 0.06551931194559302
$ws.Range("E5").Value = 0.1104196809051956
$ws.Range("G5").Value = 0.3076468350770156
$ws.Range("H5").Value = 0.473445046195728
$ws.Range("M5").Value = 0.2670308561785646
$ws.Range("O5").Value = 1.490740212568824
$ws.Range("B6").Value = 0.6208448999108782
$ws.Range("C6").Value = 0.1705253515035565
$ws.Range("D6").Value = 0.0652114979318128
$ws.Range("E6").Value = 0.1102367001066895
$ws.Range("G6").Value = 0.3076073784722055
$ws.Range("H6").Value = 0.4736446054088361
$ws.Range("M6").Value = 0.2654518484542194
$ws.Range("O6").Value = 1.491059575061428
$ws.Range("B7").Value = 0.6502788439597964
$ws.Range("C7").Value = 0.1770794723125277
$ws.Range("D7").Value = 0.06735112991303538
$ws.Range("E7").Value = 0.1115185212530605
$ws.Range("G7").Value = 0.3079156816245288
$ws.Range("H7").Value = 0.4722898010438215
$ws.Range("M7").Value = 0.2764146820809188
$ws.Range("O7").Value = 1.488978384634066
$ws.Range("B8").Value = 0.7797836566487035
$ws.Range("C8").Value = 0.2058593410120295
$ws.Range("D8").Value = 0.0768688507047699
$ws.Range("E8").Value = 0.1174675644107452
$ws.Range("G8").Value = 0.3101317229486824
$ws.Range("H8").Value = 0.4670634351289493
$ws.Range("M8").Value = 0.3248598979788184
$ws.Range("O8").Value = 1.483157973627584
$ws.Range("B9").Value = 1.032295798945086
$ws.Range("C9").Value = 0.2617646604939523
$ws.Range("D9").Value = 0.09580821940433282
$ws.Range("E9").Value = 0.1302080186562691
$ws.Range("G9").Value = 0.3175873270954099
$ws.Range("H9").Value = 0.4595271652074189
$ws.Range("M9").Value = 0.4201072830179697
$ws.Range("O9").Value = 1.483927238769212
$ws.Range("B10").Value = 1.216889297979321
$ws.Range("C10").Value = 0.3025002069578591
$ws.Range("D10").Value = 0.1098944840981062
$ws.Range("E10").Value = 0.140245773966825
$ws.Range("G10").Value = 0.3249947248658458
$ws.Range("H10").Value = 0.4556532407754332
$ws.Range("M10").Value = 0.49024534142292
$ws.Range("O10").Value = 1.492011242619725
$ws.Range("B11").Value = 1.30065611261108
$ws.Range("C11").Value = 0.320955844831218
$ws.Range("D11").Value = 0.116340825441597
$ws.Range("E11").Value = 0.1449641689515673
$ws.Range("G11").Value = 0.3287908866217464
$ws.Range("H11").Value = 0.4542540172824374
$ws.Range("M11").Value = 0.522190228723602
$ws.Range("O11").Value = 1.497341357359858
$ws.Range("B12").Value = 1.332345770163272
$ws.Range("C12").Value = 0.3279334060695476
$ws.Range("D12").Value = 0.1187874519214915
$ws.Range("E12").Value = 0.1467731584642422
$ws.Range("G12").Value = 0.3302902898558386
$ws.Range("H12").Value = 0.4537765332737536
$ws.Range("M12").Value = 0.5342925485855545
$ws.Range("O12").Value = 1.49959892564803
$ws.Range("B13").Value = 1.32552223762508
$ws.Range("C13").Value = 0.3264311670004076
$ws.Range("D13").Value = 0.1182602807069344
$ws.Range("E13").Value = 0.1463825656186231
$ws.Range("G13").Value = 0.3299646057245695
$ws.Range("H13").Value = 0.4538770361127575
$ws.Range("M13").Value = 0.5316858548596315
$ws.Range("O13").Value = 1.499102056651111
$ws.Range("B14").Value = 1.30326386961832
$ws.Range("C14").Value = 0.3215301193792186
$ws.Range("D14").Value = 0.1165419998455093
$ws.Range("E14").Value = 0.1451125479391706
$ws.Range("G14").Value = 0.3289130001386837
$ws.Range("H14").Value = 0.4542136838297921
$ws.Range("M14").Value = 0.5231857844210026
$ws.Range("O14").Value = 1.497522286862932
$ws.Range("B15").Value = 1.28962589476896
$ws.Range("C15").Value = 0.3185266179777386
$ws.Range("D15").Value = 0.1154902245917242
$ws.Range("E15").Value = 0.1443375317108533
$ws.Range("G15").Value = 0.3282769361128288
$ws.Range("H15").Value = 0.4544267155320085
$ws.Range("M15").Value = 0.5179799576082758
$ws.Range("O15").Value = 1.496585824193005
$ws.Range("B16").Value = 1.211410677173376
$ws.Range("C16").Value = 0.3012925425214519
$ws.Range("D16").Value = 0.1094739747558862
$ws.Range("E16").Value = 0.139940508031323
$ws.Range("G16").Value = 0.3247552660016169
$ws.Range("H16").Value = 0.4557520029130018
$ws.Range("M16").Value = 0.4881584399111176
$ws.Range("O16").Value = 1.491696292780091
$ws.Range("B17").Value = 1.163374462217064
$ws.Range("C17").Value = 0.2907004718728956
$ws.Range("D17").Value = 0.1057930507467404
$ws.Range("E17").Value = 0.1372823022252447
$ws.Range("G17").Value = 0.3227044634159313
$ws.Range("H17").Value = 0.4566581264544851
$ws.Range("M17").Value = 0.4698737629306038
$ws.Range("O17").Value = 1.489121091760012
$ws.Range("B18").Value = 1.135726025877204
$ws.Range("C18").Value = 0.2846011403093485
$ws.Range("D18").Value = 0.1036794970919885
$ws.Range("E18").Value = 0.1357676742720031
$ws.Range("G18").Value = 0.3215650074361065
$ws.Range("H18").Value = 0.4572134662493141
$ws.Range("M18").Value = 0.4593605508830336
$ws.Range("O18").Value = 1.487795369952408
$ws.Range("B19").Value = 1.126361473847567
$ws.Range("C19").Value = 0.2825348117858937
$ws.Range("D19").Value = 0.1029645049370487
$ws.Range("E19").Value = 0.1352572928389222
$ws.Range("G19").Value = 0.3211860796381529
$ws.Range("H19").Value = 0.4574073569666837
$ws.Range("M19").Value = 0.4558015875814618
$ws.Range("O19").Value = 1.487373158373344
$ws.Range("B20").Value = 1.168490004776402
$ws.Range("C20").Value = 0.2918287485663882
$ws.Range("D20").Value = 0.1061845170202673
$ws.Range("E20").Value = 0.1375637900377455
$ws.Range("G20").Value = 0.322918619486444
$ws.Range("H20").Value = 0.4565581312631366
$ws.Range("M20").Value = 0.471819820828415
$ws.Range("O20").Value = 1.489379125166721
$ws.Range("B21").Value = 1.309802545028276
$ws.Range("C21").Value = 0.3229699825141665
$ws.Range("D21").Value = 0.1170465502310947
$ws.Range("E21").Value = 0.1454849766199828
$ws.Range("G21").Value = 0.3292201984404102
$ws.Range("H21").Value = 0.4541133795610079
$ws.Range("M21").Value = 0.5256823129686126
$ws.Range("O21").Value = 1.497979800466652
$ws.Range("B22").Value = 1.401977012027203
$ws.Range("C22").Value = 0.3432571927738479
$ws.Range("D22").Value = 0.1241777950829572
$ws.Range("E22").Value = 0.1507917018639162
$ws.Range("G22").Value = 0.3336995772824878
$ws.Range("H22").Value = 0.4528209316045633
$ws.Range("M22").Value = 0.5609164957366914
$ws.Range("O22").Value = 1.504995603996122
$ws.Range("B23").Value = 1.352798837487285
$ws.Range("C23").Value = 0.3324356337388963
$ws.Range("D23").Value = 0.1203687585090591
$ws.Range("E23").Value = 0.1479474113678378
$ws.Range("G23").Value = 0.331275642155461
$ws.Range("H23").Value = 0.4534827420036862
$ws.Range("M23").Value = 0.5421084546879342
$ws.Range("O23").Value = 1.501123015704081
$ws.Range("B24").Value = 1.166177366861518
$ws.Range("C24").Value = 0.2913186852012757
$ws.Range("D24").Value = 0.1060075268395479
$ws.Range("E24").Value = 0.1374364870327298
$ws.Range("G24").Value = 0.3228216762859546
$ws.Range("H24").Value = 0.4566032319531104
$ws.Range("M24").Value = 0.4709400115033731
$ws.Range("O24").Value = 1.489261986193412
$ws.Range("B25").Value = 0.964144183977794
$ws.Range("C25").Value = 0.246699184481173
$ws.Range("D25").Value = 0.09065477099750296
$ws.Range("E25").Value = 0.1266440585526354
$ws.Range("G25").Value = 0.3152339408642746
$ws.Range("H25").Value = 0.4612745609518925
$ws.Range("M25").Value = 0.3943129858034808
$ws.Range("O25").Value = 1.482405357270522

Write-Host "Applied 192 cell updates for Case_4_43 (380 kV case)"
